# Adjust new architecture diagram.
#
# Nudges the icon + "NAT gateway" label for each of the three repeated
# "Private subnet" groups a few pixels, and renames the "Volume" / "Queue"
# labels to be more descriptive.

function Get-ShapeById($slide, $id) {
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $sh = $slide.Shapes.Item($i)
        if ($sh.Id -eq $id) { return $sh }
    }
    return $null
}

# PowerPoint COM exposes Shape.Left/.Top as points (single-precision floats),
# while the underlying OOXML stores EMU (1 pt = 12700 EMU). Converting
# naively (emu/12700.0) can land one EMU short once the value is rounded to
# float32 and re-quantized on save. Binary-search for a point value whose
# round-trip through float32 -> EMU reproduces the exact target EMU.
function ConvertTo-PtFromEmu($emu) {
    $lo = [double]$emu / 12700.0
    $hi = ([double]$emu + 1.0) / 12700.0
    for ($iter = 0; $iter -lt 60; $iter++) {
        $mid = ($lo + $hi) / 2.0
        $f32 = [float]$mid
        $back = [math]::Floor([double]$f32 * 12700.0)
        if ($back -lt $emu) {
            $lo = $mid
        } else {
            $hi = $mid
        }
    }
    return $hi
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Position tweaks (target EMU offsets from the OOXML diff) --------------
$moves = @(
    @{Id = 21;  X = 2671160; Y = 1485930},
    @{Id = 50;  X = 2276371; Y = 1948387},
    @{Id = 87;  X = 5898056; Y = 1507986},
    @{Id = 90;  X = 5503267; Y = 1970443},
    @{Id = 101; X = 8416197; Y = 1485930},
    @{Id = 104; X = 8021408; Y = 1948387}
)

foreach ($m in $moves) {
    $shape = Get-ShapeById $s $m.Id
    $shape.Left = ConvertTo-PtFromEmu $m.X
    $shape.Top  = ConvertTo-PtFromEmu $m.Y
}

# --- Text tweaks -------------------------------------------------------
(Get-ShapeById $s 55).TextFrame.TextRange.Text  = "EBS Volume"
(Get-ShapeById $s 61).TextFrame.TextRange.Text  = "Migration queue"
(Get-ShapeById $s 91).TextFrame.TextRange.Text  = "EBS Volume"
(Get-ShapeById $s 105).TextFrame.TextRange.Text = "EBS Volume"
